$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set the new selection (active cell) on the sheet view
$ws.Range("E15").Select()

# Rows in column B that need a value of 1 added
$rows = @(17, 21, 33, 34, 35, 36, 37, 40, 41, 42, 43, 44, 45, 46, 47)
foreach ($r in $rows) {
    $ws.Range("B$r").Value = 1
}
